# correction test assets: removed trailing space from _customVoc sheet
#
# One sheet is named "_customVoc " (note the trailing space). Find it
# (trim-matching, so this still works regardless of its position) and
# rename it to "_customVoc" (no trailing space).
$wb = $excel.ActiveWorkbook
$ws = $null
foreach ($sheet in $wb.Worksheets) {
    if ($sheet.Name.Trim() -eq "_customVoc") {
        $ws = $sheet
    }
}
$ws.Name = "_customVoc"

# That sheet's view also records a different "active" selected cell for
# the frozen bottom-left pane (B22 -> B26). Re-select the cell while the
# sheet is active so the saved sheetView reflects the new selection.
$ws.Activate()
$ws.Range("B26").Select()
